# Add two new Instagram-leader rows (155 and 156) to the bottom of the
# scraped-leads sheet, matching the "auth.json file added instead of
# login again and again" scrape re-run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 155: Arian Hushyar Van de Carr (arianvdc) ---------------------
$ws.Range("A155").Value = "Arian Hushyar Van de Carr"
$ws.Range("B155").Value = "https://www.instagram.com/arianvdc/"
$ws.Range("C155").Value = 960
$ws.Range("D155").Value = 460
$ws.Range("E155").Value = 1493
# F155 (Phone/Email) has no data for this lead, same as the rest of the sheet.
$ws.Range("G155").Value = "supplefi.co"
$ws.Range("H155").Value = "Arian Hushyar Van de Carr`narianvdc`n🌟 Jesus`n✨ Co-founder & CEO @supplefi_official`n🏡 livin’ in Austin`n❤️ married to @kedricv`n👦🏻 👧🏻 👦🏻 my babies: Knox,... `nmore`nsupplefi.co"
$ws.Range("I155").Value = "arianvdc`nFollow`n960 posts`n460 followers`n1,493 following`nArian Hushyar Van de Carr`narianvdc`n🌟 Jesus`n✨ Co-founder & CEO @supplefi_official`n🏡 livin’ in Austin`n❤️ married to @kedricv`n👦🏻 👧🏻 👦🏻 my babies: Knox,... `nmore`nsupplefi.co`nMy 40th`nHealth & Fitnes`nMy babies`nMaui"
# Match the sheet's normal (unstyled) body-row formatting used elsewhere.
$ws.Range("A155:I155").HorizontalAlignment = 1

# --- Row 156: Steve Martocci (smart) ------------------------------------
$ws.Range("A156").Value = "Steve Martocci"
$ws.Range("B156").Value = "https://www.instagram.com/smart/"
$ws.Range("C156").Value = 463
$ws.Range("D156").Value = 27100
$ws.Range("E156").Value = 1745
# F156 (Phone/Email) has no data for this lead, same as the rest of the sheet.
$ws.Range("G156").Value = "supp.co/about/founder-story"
$ws.Range("H156").Value = "Steve Martocci`nsmart`nCo-Founder of @joinsuppco, @splice , @flyblade & @groupme. Mostly Harmless.`nsupp.co/about/founder-story"
$ws.Range("I156").Value = "smart`nFollow`n463 posts`n27.1K followers`n1,745 following`nSteve Martocci`nsmart`nCo-Founder of @joinsuppco, @splice , @flyblade & @groupme. Mostly Harmless.`nsupp.co/about/founder-story`nPorto Marina`nPalisades Fire`nSuppCo`nJaxson`nSteve 4.0`nChristina`nCal`nSplice`nKomet`nChristmas '20`nSummer Tour ‘20`nThe Phish`nThanksgiving 19`nHalloween '19"
$ws.Range("A156:I156").HorizontalAlignment = 1
